$d = $word.ActiveDocument

# --- Paragraph 2: Fa5/0 (Switch) 192.168.1.1 255.255.255.0  -> ...255.255.0.0 (trailing space) ---
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("255.255.255.0", $true, $false, $false, $false, $false, $true, 0, $false, `
    "255.255.0.0 ", 2) | Out-Null

# --- Paragraph 3: S0/0/0 192.168.2.1 255.255.255.0 -> Se9/0 (Inter Router to M&B) 192.169.2.1 255.255.255.0 ---
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("S0/0/0 192.168.2.1 255.255.255.0", $true, $false, $false, $false, $false, $true, 0, $false, `
    "Se9/0 (Inter Router to M&B) 192.169.2.1 255.255.255.0", 2) | Out-Null

# --- Append the new paragraphs describing the M&B / Proctor Residence hardware setup ---
$newParagraphs = @(
    "M&B Hardware: (Done)",
    "Fa5/0 (Switch) 192.168.3.1 255.255.0.0 (Done)",
    "S0/0/0 (Inter Router to Ralph) 192.169.2.2 255.255.255.0",
    "Proctor Residence:",
    "Fa5/0 (Switch) 192.168.4.1 255.255.0.0",
    "S0/0/0 (Inter Router) 192.168.2.3 255.255.255.0",
    "",
    ""
)

$lastPara = $d.Paragraphs.Last
foreach ($text in $newParagraphs) {
    $lastPara.Range.InsertParagraphAfter()
    $lastPara = $d.Paragraphs.Last
    if ($text -ne "") {
        $lastPara.Range.Text = $text
    }
}
